$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Meerkat staging table gained a new "BusinessKey" column which the
# regenerated template places first. Shift the existing headers
# (Code, Description, Name, RoleID) one column to the right and put the
# new "BusinessKey" header in column A.
$ws.Range("E2").Value = $ws.Range("D2").Value()
$ws.Range("D2").Value = $ws.Range("C2").Value()
$ws.Range("C2").Value = $ws.Range("B2").Value()
$ws.Range("B2").Value = $ws.Range("A2").Value()
$ws.Range("A2").Value = "BusinessKey"
